# Listas sem duplicação de professores
# Clear teacher lists in rows 18-21 (columns B:F), replacing them with "-"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18:F21").Value = "-"
